$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals data (filtering save games changes the underlying
# lookup tables for TB, d2S, K and IP; G ("sum") is recomputed as
# B + C + D + E for each row).
$data = @(
    @(2, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(3, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(4, 0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059),
    @(5, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(6, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(7, 3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464),
    @(8, 0.01253208636536152, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 2.89400026249618),
    @(9, 1.445647641019636, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 2.433531715253719),
    @(10, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(11, 1.445647641019636, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 21.00072223795405),
    @(12, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(13, 0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1.104883657715537),
    @(14, 3.272327238179451, 2919.202174992006, 18.71679738969934, 13.86384647080068, 2955.055146090685),
    @(15, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(16, 0.1169995834814548, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 16.32892827181126),
    @(17, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(18, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(19, 0.6545652718822623, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 16.86649396021207),
    @(20, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(21, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(22, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(23, 1.445647641019636, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 2.433531715253719),
    @(24, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(25, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(26, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(27, 0.01253208636536152, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 0.7365606309410384),
    @(28, 0.04172184405617529, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 19.59679644099059),
    @(29, 3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464),
    @(30, 1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433),
    @(31, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(32, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(33, 3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248),
    @(34, 0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059),
    @(35, 1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455),
    @(36, 3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027),
    @(37, 1.445647641019636, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 22.32281868886277)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

Write-Output "Updated $($data.Count) rows of s_vals data"
